$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: I2 False -> True
$ws.Range("I2").Value = "'True"

# Row 3: G3 False -> True
$ws.Range("G3").Value = "'True"

# Row 4: clear all contents A4:I4 (keep formatting/style)
$ws.Range("A4:I4").ClearContents()

# Update active selection to G4
$ws.Range("G4").Select()
